$d = $word.ActiveDocument
$global:pos = 0

function ReplaceNext($old, $new) {
    $t = $d.Content.Text
    $idx = $t.IndexOf($old, $global:pos)
    if ($idx -lt 0) { throw "Text not found: $old" }
    $rng = $d.Range($idx, $idx + $old.Length)
    $rng.Text = $new
    $global:pos = $idx + $new.Length
}

function InsertAtPos($text) {
    $rng = $d.Range($global:pos, $global:pos + 1)
    $rng.Collapse(1)
    $rng.InsertBefore($text)
    $global:pos = $global:pos + $text.Length
}

ReplaceNext "type erasure and reification" "Type erasure and reification"
$global:pos = $global:pos + 1
InsertAtPos "December 05, 2018 at 05:10 Tags Programming , C & C++ , Python"
InsertAtPos " "
ReplaceNext "in this post i'd like to discuss the concepts of type erasure and reification in programming languages. i don't intend to dive very deeply into the specific rules of any particular language; rather, the post is going to present several simple examples in multiple languages, hoping to provide enough intuition and background for a more serious study, if necessary. as you'll see, the actual concepts are very simple and familiar. deeper details of specific languages pertain more to the idiosyncrasies of those languages' semantics and implementations." "In this post I'd like to discuss the concepts of type erasure and reification in programming languages. I don't intend to dive very deeply into the specific rules of any particular language; rather, the post is going to present several simple examples in multiple languages, hoping to provide enough intuition and background for a more serious study, if necessary. As you'll see, the actual concepts are very simple and familiar. Deeper details of specific languages pertain more to the idiosyncrasies of those languages' semantics and implementations."
$global:pos = $global:pos + 1
ReplaceNext "important note: in c++ there is a programming pattern called type erasure , which is quite distinct from what i'm trying to describe here [1] . i'll be using c++ examples here, but that's to demonstrate how the original concepts apply in c++. the programming pattern will be covered in a separate post." "Important note: in C++ there is a programming pattern called type erasure, which is quite distinct from what I'm trying to describe here [1]. I'll be using C++ examples here, but that's to demonstrate how the original concepts apply in C++. The programming pattern will be covered in a separate post."
$global:pos = $global:pos + 1
ReplaceNext "types at compile time, no types at run-time" "Types at compile time, no types at run-time"
$global:pos = $global:pos + 1
ReplaceNext "the title of this section is a ""one short sentence"" explanation of what type erasure means. with few exceptions, it only applies to languages with some degree of compile time (a.k.a. static ) type checking. the basic principle should be immediately familiar to folks who have some idea of what machine code generated from low-level languages like c looks like. while c has static typing, this only matters in the compiler - the generated code is completely oblivious to types." "The title of this section is a ""one short sentence"" explanation of what type erasure means. With few exceptions, it only applies to languages with some degree of compile time (a.k.a. static) type checking. The basic principle should be immediately familiar to folks who have some idea of what machine code generated from low-level languages like C looks like. While C has static typing, this only matters in the compiler - the generated code is completely oblivious to types."
$global:pos = $global:pos + 1
ReplaceNext "for example, consider the following c snippet:" "For example, consider the following C snippet:"
$global:pos = $global:pos + 1
ReplaceNext "when compiling the function extract , the compiler will perform type checking. it won't let us access fields that were not declared in the struct, for example. neither will it let us pass a pointer to a different struct (or to a float ) into extract . but once it's done helping us, the compiler generates code which is completely type-free:" "When compiling the function extract, the compiler will perform type checking. It won't let us access fields that were not declared in the struct, for example. Neither will it let us pass a pointer to a different struct (or to a float) into extract. But once it's done helping us, the compiler generates code which is completely type-free:"
$global:pos = $global:pos + 1
ReplaceNext "the compiler is familiar with the stack frame layout and other specifics of the abi, and generates code that assumes a correct type of structure was passed in. if the actual type is not what this function expects, there will be trouble (either accessing unmapped memory, or accessing wrong data)." "The compiler is familiar with the stack frame layout and other specifics of the ABI, and generates code that assumes a correct type of structure was passed in. If the actual type is not what this function expects, there will be trouble (either accessing unmapped memory, or accessing wrong data)."
$global:pos = $global:pos + 1
ReplaceNext "a slightly adjusted example will clarify this:" "A slightly adjusted example will clarify this:"
$global:pos = $global:pos + 1
ReplaceNext "the compiler will generate exactly identical code from this function, which in itself a good indication of when the types matter and when they don't. what's more interesting is that extract_cast makes it extremely easy for programmers to shoot themselves in the foot:" "The compiler will generate exactly identical code from this function, which in itself a good indication of when the types matter and when they don't. What's more interesting is that extract_cast makes it extremely easy for programmers to shoot themselves in the foot:"
$global:pos = $global:pos + 1
ReplaceNext "in general, type erasure is a concept that descibes these semantics of a language. types matter to the compiler, which uses them to generate code and help the programmer avoid errors. once everything is type-checked, however, the types are simply erased and the code the compiler generates is oblivious to them. the next section will put this in context by comparing to the opposite approach." "In general, type erasure is a concept that descibes these semantics of a language. Types matter to the compiler, which uses them to generate code and help the programmer avoid errors. Once everything is type-checked, however, the types are simply erased and the code the compiler generates is oblivious to them. The next section will put this in context by comparing to the opposite approach."
$global:pos = $global:pos + 1
ReplaceNext "while erasure means the compiler discards all type information for the actual generated code, reification is the other way to go - types are retained at run-time and used for perform various checks. a classical example from java will help demonstrate this:" "While erasure means the compiler discards all type information for the actual generated code, reification is the other way to go - types are retained at run-time and used for perform various checks. A classical example from Java will help demonstrate this:"
$global:pos = $global:pos + 1
ReplaceNext "this code creates an array of string , and converts it to a generic array of object . this is valid because arrays in java are covariant , so the compiler doesn't complain. however, in the next line we try to assign an integer into the array. this happens to fail with an exception at run-time :" "This code creates an array of String, and converts it to a generic array of Object. This is valid because arrays in Java are covariant, so the compiler doesn't complain. However, in the next line we try to assign an integer into the array. This happens to fail with an exception at run-time:"
$global:pos = $global:pos + 1
ReplaceNext "exception in thread ""main"" java.lang.arraystoreexception: java.lang.integer at main.main(main.java:5)" "Exception in thread ""main"" java.lang.ArrayStoreException: java.lang.Integer"
$global:pos = $global:pos + 1
InsertAtPos "at Main.main(Main.java:5)"
InsertAtPos " "
ReplaceNext "a type check was inserted into the generated code, and it fired when an incorrect assignment was attempted. in other words, the type of objects is reified . reification is defined roughly as ""taking something abstract and making it real/concrete"", which when applied to types means ""compile-time types are converted to actual run-time entities""." "A type check was inserted into the generated code, and it fired when an incorrect assignment was attempted. In other words, the type of objects is reified. Reification is defined roughly as ""taking something abstract and making it real/concrete"", which when applied to types means ""compile-time types are converted to actual run-time entities""."
$global:pos = $global:pos + 1
InsertAtPos "C++ has some type reification support as well, e.g. with dynamic_cast:"
InsertAtPos " "
InsertAtPos "We can call call_derived thus:"
InsertAtPos " "
ReplaceNext "the first call will successfully invoke derivedfunc ; the second will not, because the dynamic_cast will return nullptr at run-time. this is because we're using c++'s run-time type information (rtti) capabilities here, where an actual representation of the type is stored in the generated code (most likely attached to the vtable which every polymorphic object points to). c++ also has the typeid feature, but i'm showing dynamic_cast since it's the one most commonly used." "The first call will successfully invoke derivedfunc; the second will not, because the dynamic_cast will return nullptr at run-time. This is because we're using C++'s run-time type information (RTTI) capabilities here, where an actual representation of the type is stored in the generated code (most likely attached to the vtable which every polymorphic object points to). C++ also has the typeid feature, but I'm showing dynamic_cast since it's the one most commonly used."
$global:pos = $global:pos + 1
ReplaceNext "note particularly the differences between this sample and the c sample in the beginning of the post. conceptually, it's similar - we use a pointer to a general type (in c that's void* , in the c++ example we use a base type) to interact with concrete types. whereas in c there is no built-in run-time type feature, in c++ we can use rtti in some cases. with rtti enabled, dynamic_cast can be used to interact with the run-time (reified) representation of types in a limited but useful way." "Note particularly the differences between this sample and the C sample in the beginning of the post. Conceptually, it's similar - we use a pointer to a general type (in C that's void*, in the C++ example we use a base type) to interact with concrete types. Whereas in C there is no built-in run-time type feature, in C++ we can use RTTI in some cases. With RTTI enabled, dynamic_cast can be used to interact with the run-time (reified) representation of types in a limited but useful way."
$global:pos = $global:pos + 1
ReplaceNext "type erasure and java generics" "Type erasure and Java generics"
$global:pos = $global:pos + 1
ReplaceNext "one place where folks not necessarily familiar with programming language type theory encounter erasure is java generics, which were bolted onto the language after a large amount of code has already been written. the designers of java faced the binary compatibility challenge, wherein they wanted code compiled with newer java compilers to run on older vms." "One place where folks not necessarily familiar with programming language type theory encounter erasure is Java generics, which were bolted onto the language after a large amount of code has already been written. The designers of Java faced the binary compatibility challenge, wherein they wanted code compiled with newer Java compilers to run on older VMs."
$global:pos = $global:pos + 1
ReplaceNext "the solution was to use type erasure to implement generics entirely in the compiler. here's a quote from the official java generics tutorial :" "The solution was to use type erasure to implement generics entirely in the compiler. Here's a quote from the official Java generics tutorial:"
$global:pos = $global:pos + 1
ReplaceNext "generics were introduced to the java language to provide tighter type checks at compile time and to support generic programming. to implement generics, the java compiler applies type erasure to:" "Generics were introduced to the Java language to provide tighter type checks at compile time and to support generic programming. To implement generics, the Java compiler applies type erasure to:"
$global:pos = $global:pos + 1
ReplaceNext "replace all type parameters in generic types with their bounds or object if the type parameters are unbounded. the produced bytecode, therefore, contains only ordinary classes, interfaces, and methods." "Replace all type parameters in generic types with their bounds or Object if the type parameters are unbounded. The produced bytecode, therefore, contains only ordinary classes, interfaces, and methods."
$global:pos = $global:pos + 1
InsertAtPos "Insert type casts if necessary to preserve type safety."
InsertAtPos " "
ReplaceNext "generate bridge methods to preserve polymorphism in extended generic types." "Generate bridge methods to preserve polymorphism in extended generic types."
$global:pos = $global:pos + 1
ReplaceNext "here's a very simple example to demonstrate what's going on, taken from a stack overflow answer . this code:" "Here's a very simple example to demonstrate what's going on, taken from a Stack Overflow answer. This code:"
$global:pos = $global:pos + 1
ReplaceNext "uses a generic list . however, what the compiler creates prior to emitting bytecode is equivalent to:" "Uses a generic List. However, what the compiler creates prior to emitting bytecode is equivalent to:"
$global:pos = $global:pos + 1
ReplaceNext "here list is a container of object , so we can assign any element to it (similarly to the reification example shown in the previous section). the compiler then inserts a cast when accessing that element as a string. in this case the compiler will adamantly preserve type safety and won't let us do list.add(5) in the original snippet, because it sees that list is a list . therefore, the cast to (string) should be safe." "Here List is a container of Object, so we can assign any element to it (similarly to the reification example shown in the previous section). The compiler then inserts a cast when accessing that element as a string. In this case the compiler will adamantly preserve type safety and won't let us do list.add(5) in the original snippet, because it sees that list is a List. Therefore, the cast to (String) should be safe."
$global:pos = $global:pos + 1
ReplaceNext "using type erasure to implement generics with backwards compatibility is a neat idea, but it has its issues. some folks complain that not having the types available at runtime is a limitation (e.g. not being able to use instanceof and other reflection capabilities). other languages, like c# and dart 2, have reified generics which do preserve the type information at run-time." "Using type erasure to implement generics with backwards compatibility is a neat idea, but it has its issues. Some folks complain that not having the types available at runtime is a limitation (e.g. not being able to use instanceof and other reflection capabilities). Other languages, like C# and Dart 2, have reified generics which do preserve the type information at run-time."
$global:pos = $global:pos + 1
ReplaceNext "reification in dynamically typed languages" "Reification in dynamically typed languages"
$global:pos = $global:pos + 1
ReplaceNext "i hope it's obvious that the theory and techniques described above only apply to statically-typed languages. in dynamically-typed languages, like python, there is almost no concept of types at compile-time, and types are a fully reified concept. even trivial errors like:" "I hope it's obvious that the theory and techniques described above only apply to statically-typed languages. In dynamically-typed languages, like Python, there is almost no concept of types at compile-time, and types are a fully reified concept. Even trivial errors like:"
$global:pos = $global:pos + 1
ReplaceNext "fire at run-time, because there's no static type checking [2] . types obviously exist at run-time, with functions like type() and isinstance() providing complete reflection capabilities. the type() function can even create new types entirely at run-time ." "Fire at run-time, because there's no static type checking [2]. Types obviously exist at run-time, with functions like type() and isinstance() providing complete reflection capabilities. The type() function can even create new types entirely at run-time."
$global:pos = $global:pos + 1
InsertAtPos "[1] But it's most likely what you'll get to if you google for ""c++ type erasure""."
InsertAtPos " "
ReplaceNext "[2] to be clear - this is not a bug; it's a feature of python. a new method can be added to classes dynamically at runtime (here, some code could have defined a joe method for foo before the f.joe() invocation), and the compiler has absolutely no way of knowing this could or couldn't happen. so it has to assume such invocations are valid and rely on run-time checking to avoid serious errors like memory corruption." "[2] To be clear - this is not a bug; it's a feature of Python. A new method can be added to classes dynamically at runtime (here, some code could have defined a joe method for Foo before the f.joe() invocation), and the compiler has absolutely no way of knowing this could or couldn't happen. So it has to assume such invocations are valid and rely on run-time checking to avoid serious errors like memory corruption."
$global:pos = $global:pos + 1
InsertAtPos "For comments, please send me an email."
InsertAtPos " "
ReplaceNext "© 2003-2024 eli bendersky back to top" "© 2003-2024 Eli Bendersky Back to top"

Write-Host "Done. Final pos:" $global:pos
Write-Host "Final content length:" $d.Content.Text.Length
